# edit.ps1
# Applies the diff:
#   1. Inserts a new Knarot (knaerot) section - a Heading1
#      paragraph, five body paragraphs, a "Referenser - knarot"
#      Heading2 paragraph, and six reference paragraphs - right
#      after "BILAGA 1 - Fridlysta arter".
#   2. Updates the date stamp in the first-page header from
#      2023-09-13 to 2023-09-15.
#
# Strategy: first insert all the new paragraphs as blank skeletons
# (setting only their paragraph style), and only afterwards fill in
# the run text/formatting. Filling text in immediately after each
# InsertParagraphAfter() call would make the *next* paragraph mark
# inherit the character formatting (e.g. italics) of whatever run
# ends up right before the split point, leaking formatting into
# runs that should stay plain.

$d = $word.ActiveDocument

# --- Phase 0: locate the anchor paragraph ("BILAGA 1 - Fridlysta arter") ---
$baseCount = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Last

# --- Phase 1: insert blank paragraph skeletons (styles only) ---
$prev = $anchor
$prev.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($baseCount + 1)
$p1.Style = "Heading1"
$prev = $p1
$prev.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($baseCount + 2)
$p2.Style = "Normal"
$prev = $p2
$prev.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($baseCount + 3)
$p3.Style = "Normal"
$prev = $p3
$prev.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item($baseCount + 4)
$p4.Style = "Normal"
$prev = $p4
$prev.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item($baseCount + 5)
$p5.Style = "Normal"
$prev = $p5
$prev.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item($baseCount + 6)
$p6.Style = "Normal"
$prev = $p6
$prev.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item($baseCount + 7)
$p7.Style = "Heading2"
$prev = $p7
$prev.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item($baseCount + 8)
$p8.Style = "Normal"
$prev = $p8
$prev.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item($baseCount + 9)
$p9.Style = "Normal"
$prev = $p9
$prev.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Item($baseCount + 10)
$p10.Style = "Normal"
$prev = $p10
$prev.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Item($baseCount + 11)
$p11.Style = "Normal"
$prev = $p11
$prev.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Item($baseCount + 12)
$p12.Style = "Normal"
$prev = $p12
$prev.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs.Item($baseCount + 13)
$p13.Style = "Normal"
$prev = $p13

# --- Phase 2: fill in run text + formatting for each new paragraph ---
# paragraph 1/13
$rng = $p1.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Knärot – ekologi samt krav på livsmiljön")
$rng.Collapse(0) | Out-Null

# paragraph 2/13
$rng = $p2.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).")
$rng.Collapse(0) | Out-Null

# paragraph 3/13
$rng = $p3.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Samuel Johnsons doktorsavhandling ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter(" (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Vidare ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null

# paragraph 4/13
$rng = $p4.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null

# paragraph 5/13
$rng = $p5.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).")
$rng.Collapse(0) | Out-Null

# paragraph 6/13
$rng = $p6.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).")
$rng.Collapse(0) | Out-Null

# paragraph 7/13
$rng = $p7.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Referenser - knärot")
$rng.Collapse(0) | Out-Null

# paragraph 8/13
$rng = $p8.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("de Graaf M & Roberts M.R., 2009. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Short-term response of the herbaceous layer within leave patches after harvest. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Forest Ecology and Management 257, 1014-1025")
$rng.Collapse(0) | Out-Null

# paragraph 9/13
$rng = $p9.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Ecological Applications, 22, 2049-2064 ")
$rng.Collapse(0) | Out-Null

# paragraph 10/13
$rng = $p10.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Landscape Ecology, 37, sid 1839-1853")
$rng.Collapse(0) | Out-Null

# paragraph 11/13
$rng = $p11.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Biological legacies buffer local species extinction after logging. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Journal of Applied Ecology. 51, 53-62.")
$rng.Collapse(0) | Out-Null

# paragraph 12/13
$rng = $p12.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("Skogsstyrelsen, 2022. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Vägledning för hänsyn till knärot. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")
$rng.Collapse(0) | Out-Null

# paragraph 13/13
$rng = $p13.Range
$rng.Collapse(1) | Out-Null
$rng.InsertAfter("SLU Artdatabanken, 2021. ")
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Artfaktablad. Naturvård – artfakta. ")
$rng.Font.Italic = $true
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("SLU Artdatabanken, Uppsala ")
$rng.Collapse(0) | Out-Null

# --- Update the date in the "first page" header (2023-09-13 -> 2023-09-15) ---
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

Write-Output "done"
